# Append 17 months (Nov 2020 - Mar 2022) of RBNZ wholesale-rate data to the
# "Data" sheet, and bump the "Last updated" date on the "Table Description"
# sheet, matching the upstream workbook refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Data")

# --- Carry the existing number formats down onto the new rows (436:452) ---
# Row 424 is used as the style donor because it (like the target rows) has
# values in every column except K, so pasting its formats never manufactures
# a spurious K cell. Column G is only populated in some of the new rows, so
# its format is pasted solely into those.

$ws.Range("A424").Copy()
$ws.Range("A436:A452").PasteSpecial(-4122)

$ws.Range("B424:F424").Copy()
$ws.Range("B436:F452").PasteSpecial(-4122)

$ws.Range("G424").Copy()
$ws.Range("G436").PasteSpecial(-4122)
$ws.Range("G447:G452").PasteSpecial(-4122)

$ws.Range("H424").Copy()
$ws.Range("H436:H452").PasteSpecial(-4122)

$ws.Range("I424:J424").Copy()
$ws.Range("I436:J452").PasteSpecial(-4122)

$ws.Range("L424:N424").Copy()
$ws.Range("L436:N452").PasteSpecial(-4122)

$ws.Range("O424").Copy()
$ws.Range("O436:O452").PasteSpecial(-4122)

# --- Fill in the new monthly observations ---
$ws.Cells.Item(436,1).Value = 44165
$ws.Cells.Item(436,2).Value = 0.25
$ws.Cells.Item(436,3).Value = 0.22
$ws.Cells.Item(436,4).Value = 0.27
$ws.Cells.Item(436,5).Value = 0.27
$ws.Cells.Item(436,6).Value = 0.27
$ws.Cells.Item(436,7).Value = 0.21
$ws.Cells.Item(436,8).Value = 0.2
$ws.Cells.Item(436,9).Value = 0.21
$ws.Cells.Item(436,10).Value = 0.75
$ws.Cells.Item(436,12).Value = -0.53
$ws.Cells.Item(436,13).Value = -0.41
$ws.Cells.Item(436,14).Value = -0.24
$ws.Cells.Item(436,15).Value = 0.02
$ws.Cells.Item(437,1).Value = 44196
$ws.Cells.Item(437,2).Value = 0.25
$ws.Cells.Item(437,3).Value = 0.23
$ws.Cells.Item(437,4).Value = 0.26
$ws.Cells.Item(437,5).Value = 0.26
$ws.Cells.Item(437,6).Value = 0.26
$ws.Cells.Item(437,8).Value = 0.24
$ws.Cells.Item(437,9).Value = 0.34
$ws.Cells.Item(437,10).Value = 0.92
$ws.Cells.Item(437,12).Value = -0.51
$ws.Cells.Item(437,13).Value = -0.35
$ws.Cells.Item(437,14).Value = -0.09
$ws.Cells.Item(437,15).Value = 0.16
$ws.Cells.Item(438,1).Value = 44227
$ws.Cells.Item(438,2).Value = 0.25
$ws.Cells.Item(438,3).Value = 0.24
$ws.Cells.Item(438,4).Value = 0.26
$ws.Cells.Item(438,5).Value = 0.27
$ws.Cells.Item(438,6).Value = 0.28000000000000003
$ws.Cells.Item(438,8).Value = 0.22
$ws.Cells.Item(438,9).Value = 0.39
$ws.Cells.Item(438,10).Value = 1.04
$ws.Cells.Item(438,12).Value = -0.62
$ws.Cells.Item(438,13).Value = -0.42
$ws.Cells.Item(438,14).Value = -0.14000000000000001
$ws.Cells.Item(438,15).Value = 0.14000000000000001
$ws.Cells.Item(439,1).Value = 44255
$ws.Cells.Item(439,2).Value = 0.25
$ws.Cells.Item(439,3).Value = 0.24
$ws.Cells.Item(439,4).Value = 0.26
$ws.Cells.Item(439,5).Value = 0.27
$ws.Cells.Item(439,6).Value = 0.28000000000000003
$ws.Cells.Item(439,8).Value = 0.27
$ws.Cells.Item(439,9).Value = 0.75
$ws.Cells.Item(439,10).Value = 1.46
$ws.Cells.Item(439,12).Value = -0.73
$ws.Cells.Item(439,13).Value = -0.31
$ws.Cells.Item(439,14).Value = 0.13
$ws.Cells.Item(439,15).Value = 0.56999999999999995
$ws.Cells.Item(440,1).Value = 44286
$ws.Cells.Item(440,2).Value = 0.25
$ws.Cells.Item(440,3).Value = 0.23
$ws.Cells.Item(440,4).Value = 0.26
$ws.Cells.Item(440,5).Value = 0.28999999999999998
$ws.Cells.Item(440,6).Value = 0.33
$ws.Cells.Item(440,8).Value = 0.32
$ws.Cells.Item(440,9).Value = 1.03
$ws.Cells.Item(440,10).Value = 1.76
$ws.Cells.Item(440,12).Value = -0.86
$ws.Cells.Item(440,13).Value = -0.22
$ws.Cells.Item(440,14).Value = 0.33
$ws.Cells.Item(440,15).Value = 0.86
$ws.Cells.Item(441,1).Value = 44316
$ws.Cells.Item(441,2).Value = 0.25
$ws.Cells.Item(441,3).Value = 0.24
$ws.Cells.Item(441,4).Value = 0.26
$ws.Cells.Item(441,5).Value = 0.3
$ws.Cells.Item(441,6).Value = 0.34
$ws.Cells.Item(441,8).Value = 0.25
$ws.Cells.Item(441,9).Value = 0.88
$ws.Cells.Item(441,10).Value = 1.68
$ws.Cells.Item(441,12).Value = -1
$ws.Cells.Item(441,13).Value = -0.2
$ws.Cells.Item(441,14).Value = 0.44
$ws.Cells.Item(441,15).Value = 0.94
$ws.Cells.Item(442,1).Value = 44347
$ws.Cells.Item(442,2).Value = 0.25
$ws.Cells.Item(442,3).Value = 0.22
$ws.Cells.Item(442,4).Value = 0.27
$ws.Cells.Item(442,5).Value = 0.31
$ws.Cells.Item(442,6).Value = 0.35
$ws.Cells.Item(442,8).Value = 0.3
$ws.Cells.Item(442,9).Value = 1.03
$ws.Cells.Item(442,10).Value = 1.81
$ws.Cells.Item(442,12).Value = -0.99
$ws.Cells.Item(442,13).Value = -0.14000000000000001
$ws.Cells.Item(442,14).Value = 0.43
$ws.Cells.Item(442,15).Value = 0.87
$ws.Cells.Item(443,1).Value = 44377
$ws.Cells.Item(443,2).Value = 0.25
$ws.Cells.Item(443,3).Value = 0.23
$ws.Cells.Item(443,4).Value = 0.27
$ws.Cells.Item(443,5).Value = 0.3
$ws.Cells.Item(443,6).Value = 0.33
$ws.Cells.Item(443,8).Value = 0.35
$ws.Cells.Item(443,9).Value = 1.03
$ws.Cells.Item(443,10).Value = 1.76
$ws.Cells.Item(443,12).Value = -0.89
$ws.Cells.Item(443,13).Value = -0.04
$ws.Cells.Item(443,14).Value = 0.53
$ws.Cells.Item(443,15).Value = 0.97
$ws.Cells.Item(444,1).Value = 44408
$ws.Cells.Item(444,2).Value = 0.25
$ws.Cells.Item(444,3).Value = 0.23
$ws.Cells.Item(444,4).Value = 0.3
$ws.Cells.Item(444,5).Value = 0.35
$ws.Cells.Item(444,6).Value = 0.4
$ws.Cells.Item(444,8).Value = 0.65
$ws.Cells.Item(444,9).Value = 1.1299999999999999
$ws.Cells.Item(444,10).Value = 1.59
$ws.Cells.Item(444,12).Value = -0.8
$ws.Cells.Item(444,13).Value = -0.1
$ws.Cells.Item(444,14).Value = 0.33
$ws.Cells.Item(444,15).Value = 0.79
$ws.Cells.Item(445,1).Value = 44439
$ws.Cells.Item(445,2).Value = 0.25
$ws.Cells.Item(445,3).Value = 0.21
$ws.Cells.Item(445,4).Value = 0.39
$ws.Cells.Item(445,5).Value = 0.47
$ws.Cells.Item(445,6).Value = 0.54
$ws.Cells.Item(445,8).Value = 0.9
$ws.Cells.Item(445,9).Value = 1.3
$ws.Cells.Item(445,10).Value = 1.65
$ws.Cells.Item(445,12).Value = -0.72
$ws.Cells.Item(445,13).Value = -0.19
$ws.Cells.Item(445,14).Value = 0.2
$ws.Cells.Item(445,15).Value = 0.71
$ws.Cells.Item(446,1).Value = 44469
$ws.Cells.Item(446,2).Value = 0.25
$ws.Cells.Item(446,3).Value = 0.2
$ws.Cells.Item(446,4).Value = 0.38
$ws.Cells.Item(446,5).Value = 0.48
$ws.Cells.Item(446,6).Value = 0.56999999999999995
$ws.Cells.Item(446,8).Value = 1.02
$ws.Cells.Item(446,9).Value = 1.52
$ws.Cells.Item(446,10).Value = 1.87
$ws.Cells.Item(446,12).Value = -0.57999999999999996
$ws.Cells.Item(446,13).Value = [double]"-7.0000000000000007E-2"
$ws.Cells.Item(446,14).Value = 0.31
$ws.Cells.Item(446,15).Value = 0.78
$ws.Cells.Item(447,1).Value = 44500
$ws.Cells.Item(447,2).Value = 0.5
$ws.Cells.Item(447,3).Value = 0.38
$ws.Cells.Item(447,4).Value = 0.54
$ws.Cells.Item(447,5).Value = 0.62
$ws.Cells.Item(447,6).Value = 0.7
$ws.Cells.Item(447,7).Value = 1.33
$ws.Cells.Item(447,8).Value = 1.1100000000000001
$ws.Cells.Item(447,9).Value = 1.79
$ws.Cells.Item(447,10).Value = 2.21
$ws.Cells.Item(447,12).Value = -0.53
$ws.Cells.Item(447,13).Value = 0.16
$ws.Cells.Item(447,14).Value = 0.47
$ws.Cells.Item(447,15).Value = 0.78
$ws.Cells.Item(448,1).Value = 44530
$ws.Cells.Item(448,2).Value = 0.75
$ws.Cells.Item(448,3).Value = 0.54
$ws.Cells.Item(448,4).Value = 0.75
$ws.Cells.Item(448,5).Value = 0.79
$ws.Cells.Item(448,6).Value = 0.83
$ws.Cells.Item(448,7).Value = 1.58
$ws.Cells.Item(448,8).Value = 2.02
$ws.Cells.Item(448,9).Value = 2.33
$ws.Cells.Item(448,10).Value = 2.57
$ws.Cells.Item(448,12).Value = -0.35
$ws.Cells.Item(448,13).Value = 0.41
$ws.Cells.Item(448,14).Value = 0.68
$ws.Cells.Item(448,15).Value = 0.86
$ws.Cells.Item(449,1).Value = 44561
$ws.Cells.Item(449,2).Value = 0.75
$ws.Cells.Item(449,3).Value = 0.71
$ws.Cells.Item(449,4).Value = 0.79
$ws.Cells.Item(449,5).Value = 0.85
$ws.Cells.Item(449,6).Value = 0.91
$ws.Cells.Item(449,7).Value = 1.52
$ws.Cells.Item(449,8).Value = 1.97
$ws.Cells.Item(449,9).Value = 2.21
$ws.Cells.Item(449,10).Value = 2.38
$ws.Cells.Item(449,12).Value = -0.5
$ws.Cells.Item(449,13).Value = 0.23
$ws.Cells.Item(449,14).Value = 0.54
$ws.Cells.Item(449,15).Value = 0.69
$ws.Cells.Item(450,1).Value = 44592
$ws.Cells.Item(450,2).Value = 0.75
$ws.Cells.Item(450,3).Value = 0.67
$ws.Cells.Item(450,4).Value = 0.82
$ws.Cells.Item(450,5).Value = 0.92
$ws.Cells.Item(450,6).Value = 1.03
$ws.Cells.Item(450,7).Value = 1.6
$ws.Cells.Item(450,8).Value = 2.0499999999999998
$ws.Cells.Item(450,9).Value = 2.36
$ws.Cells.Item(450,10).Value = 2.56
$ws.Cells.Item(450,12).Value = -0.33
$ws.Cells.Item(450,13).Value = 0.42
$ws.Cells.Item(450,14).Value = 0.72
$ws.Cells.Item(450,15).Value = 0.85
$ws.Cells.Item(451,1).Value = 44620
$ws.Cells.Item(451,2).Value = 1
$ws.Cells.Item(451,3).Value = 0.73
$ws.Cells.Item(451,4).Value = 1
$ws.Cells.Item(451,5).Value = 1.1100000000000001
$ws.Cells.Item(451,6).Value = 1.21
$ws.Cells.Item(451,7).Value = 1.82
$ws.Cells.Item(451,8).Value = 2.25
$ws.Cells.Item(451,9).Value = 2.58
$ws.Cells.Item(451,10).Value = 2.74
$ws.Cells.Item(451,12).Value = -0.14000000000000001
$ws.Cells.Item(451,13).Value = 0.61
$ws.Cells.Item(451,14).Value = 0.87
$ws.Cells.Item(451,15).Value = 1
$ws.Cells.Item(452,1).Value = 44651
$ws.Cells.Item(452,2).Value = 1
$ws.Cells.Item(452,3).Value = 0.94
$ws.Cells.Item(452,4).Value = 1.1299999999999999
$ws.Cells.Item(452,5).Value = 1.31
$ws.Cells.Item(452,6).Value = 1.49
$ws.Cells.Item(452,7).Value = 2.14
$ws.Cells.Item(452,8).Value = 2.64
$ws.Cells.Item(452,9).Value = 2.93
$ws.Cells.Item(452,10).Value = 3.07
$ws.Cells.Item(452,12).Value = -0.17
$ws.Cells.Item(452,13).Value = 0.56999999999999995
$ws.Cells.Item(452,14).Value = 0.79
$ws.Cells.Item(452,15).Value = 0.91

# --- Bump the "last updated" date shown on the Table Description sheet ---
$ws2 = $wb.Sheets.Item("Table Description")
$ws2.Range("B4").Value = 44652
